$wb = $excel.ActiveWorkbook

# ---- Typography sheet: fill in Wildcard Characters / Wildcard Ranges ----
$wsTypo = $wb.Worksheets.Item("Typography")

foreach ($r in 4..12) {
    $wsTypo.Cells.Item($r, 7).Value = ". ,' '"   # column G = Wildcard Characters
}

foreach ($r in 10..12) {
    $wsTypo.Cells.Item($r, 9).Value = "a-z,A-Z,0-9"   # column I = Wildcard Ranges
}

# ---- Translation sheet: update display-measure labels for time mode ----
$wsTrans = $wb.Worksheets.Item("Translation")

$rowsToTrim = @(89, 92, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107)
foreach ($r in $rowsToTrim) {
    $wsTrans.Cells.Item($r, 6).Value = "IN<value>: "   # column F
}

# Append new SingleUseId143..SingleUseId156 rows (128..141) for the
# display-measure value labels of the time-mode UI.
$startRow = 128
$startId = 143
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $id = $startId + $i
    $wsTrans.Cells.Item($r, 2).Value = "SingleUseId$id"  # B = TEXT ID
    $wsTrans.Cells.Item($r, 3).Value = "displayLabel"    # C = TYPOGRAPHY NAME
    $wsTrans.Cells.Item($r, 4).Value = "Center"          # D = ALIGNMENT
    $wsTrans.Cells.Item($r, 5).Value = "LTR"             # E = DIRECTION
    $wsTrans.Cells.Item($r, 6).Value = "<value>"         # F = GB (default text)
}
